# "Update to throwing errors"
#
# The textual changes all live in the speaker notes (the Notes Placeholder
# body shape, Shapes.Item(2), on each NotesPage):
#   - slide 2 notes: prefix paragraph 1 with "(Example of CALWEB) "
#   - slide 5 notes: append a new "Show example 4" paragraph
#   - slide 6 notes: "Show example 4" -> "Show example 5"
#   - slide 7 notes: "Show example 5" -> "Show example 6"
#   - slide 8 notes: "Show example 6 then 7" -> "Show example 7 then 8"

$p = $ppt.ActivePresentation

function Set-NotesBody($slideIndex, [string[]]$paragraphs) {
    $slide = $p.Slides.Item($slideIndex)
    $notesShape = $slide.NotesPage.Shapes.Item(2)
    $notesShape.TextFrame.TextRange.Text = [string]::Join("`n", $paragraphs)
}

# --- Slide 2 notes ---
Set-NotesBody 2 @(
    '(Example of CALWEB) Environments change all the time and may have adverse affects on your scripts',
    'Someone needs to know',
    'Can the program resolve the issue itself? Could be as easy as creating a new file when one does not exist to restarting a vApp within Vcenter.'
)

# --- Slide 5 notes (new "Show example 4" paragraph added) ---
Set-NotesBody 5 @(
    'Notice I did not call this a command since this is a keyword like break and exit. There are no real cmdlets in PowerShell for these.',
    '',
    'The generic error is not a pretty one.',
    '',
    'You can clean it up by using Write-Error. The default for Write-Error is a non-terminating error',
    '',
    'Show example 3',
    'Show example 4'
)

# --- Slide 6 notes ("Show example 4" -> "Show example 5") ---
Set-NotesBody 6 @(
    'What - Catch a different error when a database connection fails versus when you divide by 0',
    '',
    'Why – When a database connection errors out your need to perform different things such as closing the database connection. Vs. dividing by zero you need to perform division slightly different. Whether that means to just default the result to 0 or null.',
    '',
    'Show example 5'
)

# --- Slide 7 notes ("Show example 5" -> "Show example 6") ---
Set-NotesBody 7 @(
    '$Error is an array of all errors in current session.',
    '',
    'If you need to perform work on a specific error you can use the ErrorVariable parameter. This is only available when you use the cmdletbinding.',
    '',
    'Show example 6'
)

# --- Slide 8 notes ("Show example 6 then 7" -> "Show example 7 then 8") ---
Set-NotesBody 8 @(
    'PowerShell’s fancy error handling only works if the function call succeeded.',
    '',
    'The alternative would be to use $PSItem',
    '',
    '$PSItem and $_ are exactly the same in regards to how they function and how fast they are. The real difference is that $_ is the v2 way and $PSItem is the v3+ way.',
    '',
    'Show example 7 then 8'
)
